$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 12:46"

# --- Swap country order: Malasia now sorts before Namibia ---
$ws.Range("A96").Value = "Malasia"
$ws.Range("A97").Value = "Namibia"

# --- Swap country order: Timor Oriental now sorts before Santa Lucia ---
$ws.Range("A207").Value = "Timor Oriental"
$ws.Range("A208").Value = "Santa Lucia"

# --- Updated daily COVID case numbers ---
$ws.Range("B4").Value = 7321465
$ws.Range("C4").Value = 122
$ws.Range("D4").Value = 4560742
$ws.Range("E4").Value = 2551269
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 209454
$ws.Range("B18").Value = 360555
$ws.Range("C18").Value = 1407
$ws.Range("D18").Value = 272073
$ws.Range("E18").Value = 83289
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = 5193
$ws.Range("B33").Value = 123944
$ws.Range("C33").Value = 1271
$ws.Range("D33").Value = 99344
$ws.Range("E33").Value = 19852
$ws.Range("G33").Value = 30
$ws.Range("H33").Value = 4748
$ws.Range("B42").Value = 98057
$ws.Range("C42").Value = 607
$ws.Range("D42").Value = 88234
$ws.Range("E42").Value = 8899
$ws.Range("G42").Value = 15
$ws.Range("H42").Value = 924
$ws.Range("B43").Value = 92095
$ws.Range("C43").Value = 626
$ws.Range("D43").Value = 81462
$ws.Range("E43").Value = 10220
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 413
$ws.Range("B61").Value = 52646
$ws.Range("C61").Value = 782
$ws.Range("E61").Value = 7882
$ws.Range("D95").Value = 7534
$ws.Range("E95").Value = 5348
$ws.Range("B96").Value = 11034
$ws.Range("C96").Value = 115
$ws.Range("D96").Value = 9889
$ws.Range("E96").Value = 1011
$ws.Range("H96").Value = 134
$ws.Range("B97").Value = 11033
$ws.Range("D97").Value = 8776
$ws.Range("E97").Value = 2137
$ws.Range("H97").Value = 120
$ws.Range("B99").Value = 10624
$ws.Range("C99").Value = 12
$ws.Range("E99").Value = 260
$ws.Range("B104").Value = 9743
$ws.Range("C104").Value = 61
$ws.Range("E104").Value = 1550
$ws.Range("B181").Value = 382
$ws.Range("C181").Value = 3
$ws.Range("E181").Value = 41
